$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.426653623580933
$ws.Range("B1").Value = 1.57176661491394
$ws.Range("C1").Value = 1.680267095565796
$ws.Range("D1").Value = 2.363569974899292
$ws.Range("E1").Value = 3.867928028106689
